# Auto-generated edit script: apply numeric corrections to Leve profit-tracking sheets
# (currentAveragePrice / LevePrice / LeveProfit columns H-N) across ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 1950
$ws.Range("I40").Value = 1950
$ws.Range("J40").Value = 1950
$ws.Range("K40").Value = 1950
$ws.Range("L40").Value = 1950
$ws.Range("M40").Value = -1775
$ws.Range("N40").Value = -2300
$ws.Range("H64").Value = 3688
$ws.Range("I64").Value = 2963.3333
$ws.Range("J64").Value = 4775
$ws.Range("K64").Value = 2963.3333
$ws.Range("L64").Value = 4775
$ws.Range("M64").Value = -2715.3333
$ws.Range("N64").Value = -5271
$ws.Range("H67").Value = 3688
$ws.Range("I67").Value = 2963.3333
$ws.Range("J67").Value = 4775
$ws.Range("K67").Value = 2963.3333
$ws.Range("L67").Value = 4775
$ws.Range("M67").Value = -2105.3333
$ws.Range("N67").Value = -6491
$ws.Range("H137").Value = 1055.8462
$ws.Range("I137").Value = 892.4737
$ws.Range("K137").Value = 2677.4211
$ws.Range("M137").Value = -127.4211

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H28").Value = 7708.4
$ws.Range("I28").Value = 6285
$ws.Range("J28").Value = 20519
$ws.Range("K28").Value = 6285
$ws.Range("L28").Value = 20519
$ws.Range("M28").Value = -6093
$ws.Range("N28").Value = -20903
$ws.Range("H74").Value = 707.25
$ws.Range("I74").Value = 689.7857
$ws.Range("J74").Value = 748
$ws.Range("K74").Value = 689.7857
$ws.Range("L74").Value = 748
$ws.Range("M74").Value = 184.2143
$ws.Range("N74").Value = -2496
$ws.Range("H77").Value = 707.25
$ws.Range("I77").Value = 689.7857
$ws.Range("J77").Value = 748
$ws.Range("K77").Value = 3448.9285
$ws.Range("L77").Value = 3740
$ws.Range("M77").Value = 919.0715
$ws.Range("N77").Value = -12476
$ws.Range("H99").Value = 7708.4
$ws.Range("I99").Value = 6285
$ws.Range("J99").Value = 20519
$ws.Range("K99").Value = 6285
$ws.Range("L99").Value = 20519
$ws.Range("M99").Value = -3290
$ws.Range("N99").Value = -26509

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H62").Value = 30000
$ws.Range("J62").Value = 30000
$ws.Range("L62").Value = 30000
$ws.Range("N62").Value = -31372
$ws.Range("H65").Value = 30000
$ws.Range("J65").Value = 30000
$ws.Range("L65").Value = 90000
$ws.Range("N65").Value = -96864
$ws.Range("H75").Value = 13000
$ws.Range("I75").Value = 13000
$ws.Range("K75").Value = 13000
$ws.Range("M75").Value = -12064
$ws.Range("H76").Value = 39314
$ws.Range("J76").Value = 39314
$ws.Range("L76").Value = 39314
$ws.Range("N76").Value = -39944
$ws.Range("H78").Value = 13000
$ws.Range("I78").Value = 13000
$ws.Range("K78").Value = 39000
$ws.Range("M78").Value = -34320
$ws.Range("H79").Value = 39314
$ws.Range("J79").Value = 39314
$ws.Range("L79").Value = 39314
$ws.Range("N79").Value = -41498
$ws.Range("H98").Value = 0
$ws.Range("J98").Value = 0
$ws.Range("L98").Value = 0
$ws.Range("N98").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2322.9375
$ws.Range("I31").Value = 1713.9166
$ws.Range("J31").Value = 4150
$ws.Range("K31").Value = 1713.9166
$ws.Range("L31").Value = 4150
$ws.Range("M31").Value = -1418.9166
$ws.Range("N31").Value = -4740
$ws.Range("H34").Value = 2322.9375
$ws.Range("I34").Value = 1713.9166
$ws.Range("J34").Value = 4150
$ws.Range("K34").Value = 1713.9166
$ws.Range("L34").Value = 4150
$ws.Range("M34").Value = -1511.9166
$ws.Range("N34").Value = -4554
$ws.Range("H97").Value = 27197
$ws.Range("J97").Value = 27197
$ws.Range("L97").Value = 27197
$ws.Range("N97").Value = -29179
$ws.Range("H99").Value = 1366.3334
$ws.Range("J99").Value = 1199.5
$ws.Range("L99").Value = 1199.5
$ws.Range("N99").Value = -4195.5
$ws.Range("H126").Value = 1366.3334
$ws.Range("J126").Value = 1199.5
$ws.Range("L126").Value = 3598.5
$ws.Range("N126").Value = -8538.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 84.3
$ws.Range("I2").Value = 17
$ws.Range("J2").Value = 129.16667
$ws.Range("K2").Value = 102
$ws.Range("L2").Value = 775.0000200000001
$ws.Range("M2").Value = 11
$ws.Range("N2").Value = -1001.00002
$ws.Range("H5").Value = 1108
$ws.Range("I5").Value = 1284
$ws.Range("J5").Value = 580
$ws.Range("K5").Value = 3852
$ws.Range("L5").Value = 1740
$ws.Range("M5").Value = -3740
$ws.Range("N5").Value = -1964
$ws.Range("H33").Value = 325.3125
$ws.Range("I33").Value = 83.181816
$ws.Range("J33").Value = 858
$ws.Range("K33").Value = 499.090896
$ws.Range("L33").Value = 5148
$ws.Range("M33").Value = -216.090896
$ws.Range("N33").Value = -5714
$ws.Range("H107").Value = 358.43478
$ws.Range("J107").Value = 356.33334
$ws.Range("L107").Value = 1069.00002
$ws.Range("N107").Value = -4909.000019999999
$ws.Range("H131").Value = 11628973
$ws.Range("J131").Value = 14706978
$ws.Range("L131").Value = 44120934
$ws.Range("N131").Value = -44131014
$ws.Range("H135").Value = 1108
$ws.Range("I135").Value = 1284
$ws.Range("J135").Value = 580
$ws.Range("K135").Value = 11556
$ws.Range("L135").Value = 5220
$ws.Range("M135").Value = -9021
$ws.Range("N135").Value = -10290

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 18152772
$ws.Range("I70").Value = 35163464
$ws.Range("J70").Value = 8033.933
$ws.Range("K70").Value = 35163464
$ws.Range("L70").Value = 8033.933
$ws.Range("M70").Value = -35163194
$ws.Range("N70").Value = -8573.933000000001
$ws.Range("H73").Value = 18152772
$ws.Range("I73").Value = 35163464
$ws.Range("J73").Value = 8033.933
$ws.Range("K73").Value = 35163464
$ws.Range("L73").Value = 8033.933
$ws.Range("M73").Value = -35162528
$ws.Range("N73").Value = -9905.933000000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H32").Value = 0
$ws.Range("I32").Value = 0
$ws.Range("K32").Value = 0
$ws.Range("M32").ClearContents()
$ws.Range("H47").Value = 5399
$ws.Range("J47").Value = 5399
$ws.Range("L47").Value = 5399
$ws.Range("N47").Value = -6379
$ws.Range("H52").Value = 5399
$ws.Range("J52").Value = 5399
$ws.Range("L52").Value = 5399
$ws.Range("N52").Value = -5865
$ws.Range("H68").Value = 1964.3077
$ws.Range("I68").Value = 1831.68
$ws.Range("J68").Value = 2201.1428
$ws.Range("K68").Value = 1831.68
$ws.Range("L68").Value = 2201.1428
$ws.Range("M68").Value = -1082.68
$ws.Range("N68").Value = -3699.1428
$ws.Range("H71").Value = 1964.3077
$ws.Range("I71").Value = 1831.68
$ws.Range("J71").Value = 2201.1428
$ws.Range("K71").Value = 9158.4
$ws.Range("L71").Value = 11005.714
$ws.Range("M71").Value = -5414.4
$ws.Range("N71").Value = -18493.714
$ws.Range("H99").Value = 0
$ws.Range("I99").Value = 0
$ws.Range("K99").Value = 0
$ws.Range("M99").ClearContents()
$ws.Range("H136").Value = 7514.423
$ws.Range("I136").Value = 6660.7144
$ws.Range("J136").Value = 11100
$ws.Range("K136").Value = 19982.1432
$ws.Range("L136").Value = 33300
$ws.Range("M136").Value = -17432.1432
$ws.Range("N136").Value = -38400

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 2389.5
$ws.Range("I122").Value = 2468.3
$ws.Range("J122").Value = 2192.5
$ws.Range("K122").Value = 7404.900000000001
$ws.Range("L122").Value = 6577.5
$ws.Range("M122").Value = -4954.900000000001
$ws.Range("N122").Value = -11477.5
